$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price/Volume columns as text so Excel does not silently
# reformat values such as "245.69" or "1.001" into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Update price (D) and volume-change (E) values for rows that only changed their figures
$ws.Range("D2").Value = "30.407.36"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.869.55"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "245.69"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D7").Value = "0.4740"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("D9").Value = "0.06491"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "22.06"
$ws.Range("E10").Value = "  +5.87%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "97.48"
$ws.Range("E12").Value = "  +2.78%  "
$ws.Range("D13").Value = "0.7359"
$ws.Range("E13").Value = "  +7.05%  "
$ws.Range("D14").Value = "1.869.89"
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").Value = "5.127"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "273.45"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").Value = "30.405.67"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "13.35"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "0.000007528"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "2.116.70"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "5.216"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "6.165"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").Value = "9.281"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").Value = "163.73"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("D27").Value = "18.79"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "1.924"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("D29").Value = "0.09994"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").Value = "1.499"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "4.294"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").Value = "4.117"
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("D34").Value = "0.04825"
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("D35").Value = "1.123"
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("D36").Value = "0.6961"
$ws.Range("E36").Value = "  +1.22%  "

# Rows 37-51: a new "Frax" entry was inserted on 2023-07-11, shifting the remaining
# coins down by one slot and pushing "Decentraland" off the bottom of the list.
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "0.9999"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "2.710"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01851"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.750"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.303"
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "73.52"
$ws.Range("E42").Value = "  +4.27%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "1.961"
$ws.Range("E43").Value = "  +3.23%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4181"
$ws.Range("E44").Value = "  +3.04%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "0.8337"
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "101.87"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.275"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "35.46"
$ws.Range("E49").Value = "  +2.83%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "6.997"
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "922.76"
$ws.Range("E51").Value = "  +0.06%  "